# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Updates the "K" column (column G) values for rows 2-15 on Sheet1 to the
# newly-calculated strikeout/K figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 6
    3  = 2
    4  = 0
    5  = 3
    6  = 5
    7  = 2
    8  = 3
    9  = 2
    10 = 1
    11 = 3
    12 = 6
    13 = 4
    14 = 2
    15 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
